$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.626.82'
$ws.Range("E2").Value = '  +0.47%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.707.09'
$ws.Range("E3").Value = '  +0.88%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '673.32'
$ws.Range("E5").Value = '  -1.26%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.02'
$ws.Range("E6").Value = '  +2.18%  '

# Row 7
$ws.Range("E7").Value = '  +0.06%  '

# Row 8
$ws.Range("E8").Value = '  +1.03%  '

# Row 9
$ws.Range("E9").Value = '  +0.75%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.10'
$ws.Range("E10").Value = '  +1.91%  '

# Row 11
$ws.Range("E11").Value = '  +1.92%  '

# Row 12
$ws.Range("E12").Value = '  +1.26%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '32.87'
$ws.Range("E13").Value = '  +2.15%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.689.32'
$ws.Range("E14").Value = '  +0.42%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '69.660.63'
$ws.Range("E15").Value = '  +0.53%  '

# Row 16
$ws.Range("E16").Value = '  +1.74%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '16.34'
$ws.Range("E17").Value = '  +2.99%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.53'
$ws.Range("E18").Value = '  +2.35%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '474.19'
$ws.Range("E19").Value = '  +1.14%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.82'
$ws.Range("E20").Value = '  -1.51%  '

# Row 21
$ws.Range("E21").Value = '  +0.98%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '80.47'
$ws.Range("E22").Value = '  +0.63%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.856.13'
$ws.Range("E23").Value = '  +0.91%  '

# Row 24
$ws.Range("E24").Value = '  +5.63%  '

# Row 25
$ws.Range("E25").Value = '  -0.02%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.02'
$ws.Range("E26").Value = '  +1.24%  '

# Row 27
$ws.Range("E27").Value = '  +0.46%  '

# Row 28
$ws.Range("E28").Value = '  -0.13%  '

# Row 29
$ws.Range("E29").Value = '  +0.00%  '

# Row 30
$ws.Range("E30").Value = '  +1.51%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.59'
$ws.Range("E31").Value = '  +0.58%  '

# Row 32
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.167'
$ws.Range("E32").Value = '  +6.27%  '

# Row 33
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.96'
$ws.Range("E33").Value = '  +0.65%  '

# Row 34
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  -0.14%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.697.40'
$ws.Range("E35").Value = '  +1.18%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.55'
$ws.Range("E36").Value = '  +4.74%  '

# Row 37
$ws.Range("E37").Value = '  +1.56%  '

# Row 38
$ws.Range("E38").Value = '  -0.02%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.26'
$ws.Range("E39").Value = '  +1.62%  '

# Row 40
$ws.Range("E40").Value = '  -0.06%  '

# Row 41
$ws.Range("E41").Value = '  +1.61%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '174.05'
$ws.Range("E42").Value = '  +4.38%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.941'
$ws.Range("E43").Value = '  +0.22%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '47.07'
$ws.Range("E44").Value = '  -1.11%  '

# Row 45
$ws.Range("E45").Value = '  +1.54%  '

# Row 46
$ws.Range("E46").Value = '  +1.19%  '

# Row 47
$ws.Range("E47").Value = '  +1.95%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.72'
$ws.Range("E48").Value = '  +2.89%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.10'
$ws.Range("E49").Value = '  -0.31%  '

# Row 50
$ws.Range("E50").Value = '  +1.84%  '

# Row 51
$ws.Range("E51").Value = '  +1.50%  '
